# Adição de tempo de preparo na receita
$wb = $excel.ActiveWorkbook

# --- Sheet "Receita": insert a new column before "tipo" (column C) ---
$wsReceita = $wb.Worksheets.Item("Receita")
$wsReceita.Columns.Item(3).Insert()
$wsReceita.Range("C1").Value = "tempoPreparo"

# --- Sheet "Ingrediente": drop the long tail of empty formatted rows ---
$wsIngrediente = $wb.Worksheets.Item("Ingrediente")
$wsIngrediente.Columns.Item(1).ClearFormats()
$wsIngrediente.Columns.Item(2).ClearFormats()
$wsIngrediente.Columns.Item(3).ClearFormats()
$wsIngrediente.Range("D1").ClearFormats()
$wsIngrediente.Range("A2:D27").EntireRow.Delete()

# --- Sheet "Categoria": drop the long tail of empty formatted rows ---
$wsCategoria = $wb.Worksheets.Item("Categoria")
$wsCategoria.Columns.Item(1).ClearFormats()
$wsCategoria.Range("A1:B1").ClearFormats()
$wsCategoria.Range("A2:B6").EntireRow.Delete()

# --- Restore the selection / active tab so "Receita" is the active sheet ---
$wsReceita.Select()
$wsReceita.Range("C2").Select()
